{"js": "// \"template updates, font control\"\n//\n// Target change (word/styles.xml): the \"Hyperlink\" character style now\n// bases itself on \"VerbatimChar\" (styleId of the \"Verbatim Char\" style)\n// instead of \"CaptionChar\", and picks up an explicit 11pt (sz=22\n// half-points) run-properties override alongside the existing single\n// underline.\n//\n// (The diff also touches the opaque w:tmpl nonces on the ten unused,\n// Word-internal \"legacy single-level list\" abstractNum defs in\n// word/numbering.xml. Those hex codes are bookkeeping Word mints for its\n// own list-gallery cache \u2014 there is no Word JS API / object-model surface\n// that reads or writes them, so that part of the diff is not reachable\n// from script and is intentionally left alone here.)\n\nconst styles = context.document.styles;\nconst hyperlink = styles.getByNameOrNullObject(\"Hyperlink\");\nhyperlink.load(\"baseStyle\");\nawait context.sync();\n\nif (!hyperlink.isNullObject) {\n  // Re-parent the style onto \"Verbatim Char\" (styleId \"VerbatimChar\").\n  hyperlink.baseStyle = \"VerbatimChar\";\n  // Add the explicit 11pt font-size override (w:sz w:val=\"22\").\n  hyperlink.font.size = 11;\n  await context.sync();\n}\n", "ps1": "# \"template updates, font control\"\n#\n# Target change (word/styles.xml): the \"Hyperlink\" character style now\n# bases itself on \"VerbatimChar\" (styleId of the \"Verbatim Char\" style)\n# instead of \"CaptionChar\", and picks up an explicit 11pt (sz=22\n# half-points) run-properties override alongside the existing single\n# underline.\n#\n# (The diff also touches the opaque w:tmpl nonces on the ten unused,\n# Word-internal \"legacy single-level list\" abstractNum defs in\n# word/numbering.xml. Those hex codes are bookkeeping Word mints for its\n# own list-gallery cache -- there is no Word COM object-model property\n# that reads or writes them (no AbstractNum/Tmpl member exists on Style,\n# ListTemplate, ListGalleries, etc.), so that part of the diff is not\n# reachable from script and is intentionally left alone here.)\n\n$d = $word.ActiveDocument\n\n$hyperlink = $d.Styles(\"Hyperlink\")\nif ($hyperlink -ne $null) {\n    # Re-parent the style onto \"Verbatim Char\" (styleId \"VerbatimChar\").\n    $hyperlink.BaseStyle = $d.Styles(\"VerbatimChar\")\n    # Add the explicit 11pt font-size override (w:sz w:val=\"22\").\n    $hyperlink.Font.Size = 11\n}\n"}
